$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values could be misread as numbers (trailing-zero
# decimals like "1.00" / "159.90") need to be forced to Text first,
# then restored to the default (Normal/General) style so the saved
# file keeps using the original unstyled cell format.
$textForceCells = @("D4", "D15", "D31", "D32", "D36", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '41.723.33'
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").Value = '2.472.45'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.29%  '

$ws.Range("D5").Value = '318.43'
$ws.Range("E5").Value = '  +1.10%  '

$ws.Range("D6").Value = '93.22'
$ws.Range("E6").Value = '  +1.54%  '

$ws.Range("E7").Value = '  +0.37%  '

$ws.Range("D9").Value = '0.516'
$ws.Range("E9").Value = '  +0.11%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").Value = '0.0869'
$ws.Range("E10").Value = '  +9.47%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").Value = '33.41'
$ws.Range("E11").Value = '  +2.93%  '

$ws.Range("E12").Value = '  +0.63%  '

$ws.Range("D13").Value = '2.854.61'
$ws.Range("E13").Value = '  +0.17%  '

$ws.Range("E14").Value = '  +0.77%  '

$ws.Range("D15").Value = '15.70'
$ws.Range("E15").Value = '  -2.17%  '

$ws.Range("D16").Value = '2.458.41'
$ws.Range("E16").Value = '  -0.38%  '

$ws.Range("D17").Value = '0.793'
$ws.Range("E17").Value = '  +2.16%  '

$ws.Range("D18").Value = '41.692.06'
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("E19").Value = '  -0.70%  '

$ws.Range("E20").Value = '  +0.51%  '

$ws.Range("D21").Value = '71.09'
$ws.Range("E21").Value = '  -0.13%  '

$ws.Range("D22").Value = '11.27'
$ws.Range("E22").Value = '  +1.58%  '

$ws.Range("D23").Value = '239.39'
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("E24").Value = '  +0.72%  '

$ws.Range("D25").Value = '1.94'
$ws.Range("E25").Value = '  +1.56%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("E27").Value = '  -0.42%  '

$ws.Range("E28").Value = '  +0.86%  '

$ws.Range("D29").Value = '9.78'
$ws.Range("E29").Value = '  +0.97%  '

$ws.Range("D30").Value = '36.01'
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("D31").Value = '159.90'
$ws.Range("E31").Value = '  +2.72%  '

$ws.Range("D32").Value = '5.50'
$ws.Range("E32").Value = '  +0.93%  '

$ws.Range("E33").Value = '  -0.13%  '

$ws.Range("E34").Value = '  +0.30%  '

$ws.Range("D35").Value = '0.0764'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("D36").Value = '17.50'
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("E37").Value = '  +4.64%  '

$ws.Range("E38").Value = '  +1.29%  '

$ws.Range("E39").Value = '  +1.71%  '

$ws.Range("D40").Value = '0.104'
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("D41").Value = '4.01'
$ws.Range("E41").Value = '  +0.32%  '

$ws.Range("D42").Value = '2.51'
$ws.Range("E42").Value = '  +6.41%  '

$ws.Range("D43").Value = '1.997.88'
$ws.Range("E43").Value = '  +2.85%  '

$ws.Range("E44").Value = '  +0.69%  '

$ws.Range("D45").Value = '18.91'
$ws.Range("E45").Value = '  +0.96%  '

$ws.Range("E46").Value = '  +2.81%  '

$ws.Range("D47").Value = '9.51'
$ws.Range("E47").Value = '  +4.67%  '

$ws.Range("D48").Value = '2.711.16'
$ws.Range("E48").Value = '  +0.12%  '

$ws.Range("D49").Value = '97.54'
$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").Value = '73.80'
$ws.Range("E50").Value = '  +3.10%  '

$ws.Range("D51").Value = '67.16'
$ws.Range("E51").Value = '  -0.08%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}